$d = $word.ActiveDocument

# --- Simple single-run text substitutions (title + body paragraph sentences) ---
$d.Content.Find.Execute('Unraveling the Enigma of Dark Matter', $true, $true, $false, $false, $false, $true, 1, $false, 'Chemistry and the Symphony of Life', 2) | Out-Null
$d.Content.Find.Execute('As we gaze upon the vast cosmos, a profound mystery beckons us: the enigmatic essence known as dark matter', $true, $true, $false, $false, $false, $true, 1, $false, 'Chemistry, the study of matter and its interactions, holds the key to unlocking the intricate secrets of life', 2) | Out-Null
$d.Content.Find.Execute(' Comprising approximately 27% of the universe, this enigmatic substance eludes our direct observation, revealing its presence only through its gravitational influence on the cosmos', $true, $true, $false, $false, $false, $true, 1, $false, ' With its vast array of elements and compounds, chemistry forms the very foundation of our natural world, orchestrating the delicate balance that sustains all living organisms', 2) | Out-Null
$d.Content.Find.Execute(' What is the true nature of dark matter? How can we unravel its secrets? These questions have ignited a relentless pursuit among scientists, leading to captivating discoveries and fueling our understanding of the universe''s composition and evolution', $true, $true, $false, $false, $false, $true, 1, $false, ' From the smallest building blocks of life to the complex processes that govern our bodies, chemistry weaves a symphony of intricate reactions that are essential for our existence', 2) | Out-Null
$d.Content.Find.Execute('The existence of dark matter first emerged as a paradox in the 1930s when astronomers observed the rotational speeds of galaxies', $true, $true, $false, $false, $false, $true, 1, $false, 'As we delve into the world of chemistry, we embark on a journey of discovery, exploring the properties and behaviors of substances, their interactions, and their role in shaping our lives', 2) | Out-Null
$d.Content.Find.Execute(' Contrary to expectations, the stars at the outskirts of galaxies exhibited velocities that defied the laws of physics, suggesting the presence of unseen mass exerting gravitational influence', $true, $true, $false, $false, $false, $true, 1, $false, ' We uncover the mysteries behind the transformation of matter, from the combustion of fuels to the formation of new materials, unraveling the intricate dance of electrons, atoms, and molecules', 2) | Out-Null
$d.Content.Find.Execute(' This discrepancy served as a crucial clue, hinting at the existence of a mysterious substance dominating the universe', $true, $true, $false, $false, $false, $true, 1, $false, ' Chemistry allows us to understand the composition of the substances that surround us, from the air we breathe to the food we eat, opening up a world of possibilities for innovation and progress', 2) | Out-Null
$d.Content.Find.Execute('As scientists delved deeper into this enigma, additional evidence emerged corroborating the existence of dark matter', $true, $true, $false, $false, $false, $true, 1, $false, 'Through chemistry, we gain insight into the intricate mechanisms that drive our bodies, from the intricate workings of our cells to the complex interactions of hormones and neurotransmitters', 2) | Out-Null
$d.Content.Find.Execute(' Gravitational lensing, the bending of light due to the presence of mass, provided compelling evidence of dark matter''s gravitational effects', $true, $true, $false, $false, $false, $true, 1, $false, ' We discover the secrets behind the remarkable transformations that occur within us, from the synthesis of proteins to the intricate process of metabolism', 2) | Out-Null
$d.Content.Find.Execute(' Furthermore, observations of galaxy clusters revealed that the mass inferred from gravitational measurements far exceeded the mass contributed by visible matter, further solidifying the case for dark matter''s existence', $true, $true, $false, $false, $false, $true, 1, $false, ' Chemistry helps us navigate the challenges of illness, offering solutions through the development of medicines and treatments that target specific ailments, alleviating human suffering and enhancing our collective well-being', 2) | Out-Null

# --- Byline: "Alex Barfield" -> "Dr" + "." + " Geraldine Williams" (3 runs) ---
$pName = $d.Paragraphs(2).Range
$pNameSub = $d.Range($pName.Start, $pName.End - 1)
$pNameSub.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Dr</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> Geraldine Williams</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Email: "abarfield@protonmail"."com" -> "geraldine"."williams@loyalschool"."org" (5 runs) ---
$pEmail = $d.Paragraphs(3).Range
$pEmailSub = $d.Range($pEmail.Start, $pEmail.End - 1)
$pEmailSub.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>geraldine</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>williams@loyalschool</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>org</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Summary body paragraph: 8 runs -> 7 runs (merge + lastRenderedPageBreak) ---
$pSum = $d.Paragraphs(7).Range
$pSumSub = $d.Range($pSum.Start, $pSum.End - 1)
$pSumSub.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>In this essay, we explored the captivating realm of chemistry, revealing its crucial role in shaping our lives and our understanding of the natural world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> From the fundamental building blocks of matter to the complex symphony of reactions that orchestrate life, chemistry stands as a testament to the boundless wonders of the universe</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> As we continue to uncover the intricacies of this subject, we unlock new possibilities for innovation, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>progress, and the betterment of humanity, ensuring a future brimming with transformative discoveries and endless possibilities</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Append a new empty paragraph at the very end of the body ---
$d.Content.InsertParagraphAfter()
